$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr1 = New-Object 'object[,]' 24,5
$arr1[0,0] = 1.019999999999999
$arr1[0,1] = 1.051327715307378
$arr1[0,2] = 1.051764338027758
$arr1[0,3] = 1.060291986558754
$arr1[0,4] = 1.068328339238921
$arr1[1,0] = 1.02
$arr1[1,1] = 1.05255557359645
$arr1[1,2] = 1.05271545701331
$arr1[1,3] = 1.061703087550109
$arr1[1,4] = 1.069846996974822
$arr1[2,0] = 1.02
$arr1[2,1] = 1.053349172494521
$arr1[2,2] = 1.053330048667331
$arr1[2,3] = 1.062616086598697
$arr1[2,4] = 1.070829675353788
$arr1[3,0] = 1.02
$arr1[3,1] = 1.053682587491913
$arr1[3,2] = 1.053588222308022
$arr1[3,3] = 1.062999897271354
$arr1[3,4] = 1.071242799432319
$arr1[4,0] = 1.02
$arr1[4,1] = 1.053738556884417
$arr1[4,2] = 1.053631559082598
$arr1[4,3] = 1.063064340073225
$arr1[4,4] = 1.071312165279998
$arr1[5,0] = 1.02
$arr1[5,1] = 1.053353628438743
$arr1[5,2] = 1.053333499182469
$arr1[5,3] = 1.062621215145642
$arr1[5,4] = 1.070835195508011
$arr1[6,0] = 1.02
$arr1[6,1] = 1.051742865652259
$arr1[6,2] = 1.052085949400182
$arr1[6,3] = 1.060768893512195
$arr1[6,4] = 1.068841578546472
$arr1[7,0] = 1.02
$arr1[7,1] = 1.048897397758108
$arr1[7,2] = 1.049881045939212
$arr1[7,3] = 1.057504071428364
$arr1[7,4] = 1.06532839487847
$arr1[8,0] = 1.02
$arr1[8,1] = 1.046995438207687
$arr1[8,2] = 1.048406580735505
$arr1[8,3] = 1.055326708047283
$arr1[8,4] = 1.062985850373711
$arr1[9,0] = 1.02
$arr1[9,1] = 1.046170640169184
$arr1[9,2] = 1.047767019569463
$arr1[9,3] = 1.054383629321915
$arr1[9,4] = 1.061971334583277
$arr1[10,0] = 1.02
$arr1[10,1] = 1.045864083448707
$arr1[10,2] = 1.04752928917426
$arr1[10,3] = 1.054033282438325
$arr1[10,4] = 1.061594465490636
$arr1[11,0] = 1.02
$arr1[11,1] = 1.04592984957102
$arr1[11,2] = 1.047580290833501
$arr1[11,3] = 1.054108435179438
$arr1[11,4] = 1.061675306768129
$arr1[12,0] = 1.02
$arr1[12,1] = 1.046145303995214
$arr1[12,2] = 1.047747372160536
$arr1[12,3] = 1.054354670498316
$arr1[12,4] = 1.061940183169133
$arr1[13,0] = 1.02
$arr1[13,1] = 1.046278027274576
$arr1[13,2] = 1.047850294026191
$arr1[13,3] = 1.054506378017311
$arr1[13,4] = 1.062103377743394
$arr1[14,0] = 1.02
$arr1[14,1] = 1.047050150920366
$arr1[14,2] = 1.048449002751823
$arr1[14,3] = 1.055389291003779
$arr1[14,4] = 1.063053176177857
$arr1[15,0] = 1.02
$arr1[15,1] = 1.047534149573498
$arr1[15,2] = 1.048824258558751
$arr1[15,3] = 1.055943043906606
$arr1[15,4] = 1.063648907657706
$arr1[16,0] = 1.02
$arr1[16,1] = 1.047816338729568
$arr1[16,2] = 1.049043032182437
$arr1[16,3] = 1.056266013312471
$arr1[16,4] = 1.063996370925068
$arr1[17,0] = 1.02
$arr1[17,1] = 1.047912537868162
$arr1[17,2] = 1.049117610302643
$arr1[17,3] = 1.056376133408695
$arr1[17,4] = 1.064114844244646
$arr1[18,0] = 1.02
$arr1[18,1] = 1.047482233445923
$arr1[18,2] = 1.048784008221276
$arr1[18,3] = 1.055883634118283
$arr1[18,4] = 1.063584993092792
$arr1[19,0] = 1.02
$arr1[19,1] = 1.046081863315926
$arr1[19,2] = 1.047698175556289
$arr1[19,3] = 1.054282161633127
$arr1[19,4] = 1.06186218462291
$arr1[20,0] = 1.02
$arr1[20,1] = 1.045200293649901
$arr1[20,2] = 1.047014491030666
$arr1[20,3] = 1.053274986598864
$arr1[20,4] = 1.060778794084027
$arr1[21,0] = 1.02
$arr1[21,1] = 1.045667736201818
$arr1[21,2] = 1.047377018695141
$arr1[21,3] = 1.053808936166202
$arr1[21,4] = 1.061353140137373
$arr1[22,0] = 1.02
$arr1[22,1] = 1.047505692476997
$arr1[22,2] = 1.048802195945991
$arr1[22,3] = 1.055910478922198
$arr1[22,4] = 1.063613873385151
$arr1[23,0] = 1.02
$arr1[23,1] = 1.049633881984684
$arr1[23,2] = 1.05045185598863
$arr1[23,3] = 1.058348229238945
$arr1[23,4] = 1.066236689118706

$arr2 = New-Object 'object[,]' 24,6
$arr2[0,0] = 1.044262139707192
$arr2[0,1] = 1.056356170754503
$arr2[0,2] = 1.054514935958826
$arr2[0,3] = 1.063019163946364
$arr2[0,4] = 1.071033815629506
$arr2[0,5] = 1.057856317682721
$arr2[1,0] = 1.044622305286532
$arr2[1,1] = 1.057232711179185
$arr2[1,2] = 1.05527841056208
$arr2[1,3] = 1.064243155258974
$arr2[1,4] = 1.072366684853348
$arr2[1,5] = 1.058734102895347
$arr2[2,0] = 1.044853713886055
$arr2[2,1] = 1.05779846427092
$arr2[2,2] = 1.055770963234255
$arr2[2,3] = 1.065034554386292
$arr2[2,4] = 1.073228637057632
$arr2[2,5] = 1.059300659421365
$arr2[3,0] = 1.044950605831975
$arr2[3,1] = 1.058035967458245
$arr2[3,2] = 1.055977683142229
$arr2[3,3] = 1.065367117729811
$arr2[3,4] = 1.07359088554677
$arr2[3,5] = 1.059538499890457
$arr2[4,0] = 1.044966851444954
$arr2[4,1] = 1.05807582545054
$arr2[4,2] = 1.056012371898219
$arr2[4,3] = 1.06542294843549
$arr2[4,4] = 1.073651701971108
$arr2[4,5] = 1.059578414485674
$arr2[5,0] = 1.044855010100913
$arr2[5,1] = 1.057801639130995
$arr2[5,2] = 1.055773726805767
$arr2[5,3] = 1.065038998662415
$arr2[5,4] = 1.073233477891007
$arr2[5,5] = 1.059303838790105
$arr2[6,0] = 1.044384200078075
$arr2[6,1] = 1.056652698562348
$arr2[6,2] = 1.054773260521028
$arr2[6,3] = 1.063432945248095
$arr2[6,4] = 1.07148437181432
$arr2[6,5] = 1.058153266594065
$arr2[7,0] = 1.043541946365811
$arr2[7,1] = 1.054617086401239
$arr2[7,2] = 1.052998997073005
$arr2[7,3] = 1.060598059765218
$arr2[7,4] = 1.068398176985611
$arr2[7,5] = 1.0561147636302
$arr2[8,0] = 1.042971889261711
$arr2[8,1] = 1.053252449985237
$arr2[8,2] = 1.051808429923765
$arr2[8,3] = 1.058704646552231
$arr2[8,4] = 1.066337723659775
$arr2[8,5] = 1.054748189273961
$arr2[9,0] = 1.042723004215075
$arr2[9,1] = 1.052659721034021
$arr2[9,2] = 1.051291042045424
$arr2[9,3] = 1.057883893385413
$arr2[9,4] = 1.065444754976793
$arr2[9,5] = 1.054154618579648
$arr2[10,0] = 1.042630248397379
$arr2[10,1] = 1.052439276894842
$arr2[10,2] = 1.051098578624993
$arr2[10,3] = 1.057578890040083
$arr2[10,4] = 1.0651129443648
$arr2[10,5] = 1.053933861384505
$arr2[11,0] = 1.042650158832222
$arr2[11,1] = 1.05248657554794
$arr2[11,2] = 1.051139875505816
$arr2[11,3] = 1.057644320666159
$arr2[11,4] = 1.06518412448603
$arr2[11,5] = 1.053981227207115
$arr2[12,0] = 1.042715343298163
$arr2[12,1] = 1.052641504745179
$arr2[12,2] = 1.051275138736508
$arr2[12,3] = 1.057858684579536
$arr2[12,4] = 1.065417329933261
$arr2[12,5] = 1.054136376421585
$arr2[13,0] = 1.042755464635303
$arr2[13,1] = 1.052736924782974
$arr2[13,2] = 1.051358441367316
$arr2[13,3] = 1.057990742684243
$arr2[13,4] = 1.065560999143975
$arr2[13,5] = 1.054231931966779
$arr2[14,0] = 1.04298836370825
$arr2[14,1] = 1.053291748598568
$arr2[14,2] = 1.051842727765864
$arr2[14,3] = 1.05875909800562
$arr2[14,4] = 1.066396970191998
$arr2[14,5] = 1.054787543695831
$arr2[15,0] = 1.043133906265659
$arr2[15,1] = 1.053639282282824
$arr2[15,2] = 1.052146006964238
$arr2[15,3] = 1.059240824496295
$arr2[15,4] = 1.066921140518898
$arr2[15,5] = 1.055135570917781
$arr2[16,0] = 1.043218601406364
$arr2[16,1] = 1.053841816202248
$arr2[16,2] = 1.052322724862228
$arr2[16,3] = 1.059521721606766
$arr2[16,4] = 1.067226805418346
$arr2[16,5] = 1.055338392458602
$arr2[17,0] = 1.043247446823887
$arr2[17,1] = 1.05391084516837
$arr2[17,2] = 1.05238295068622
$arr2[17,3] = 1.059617485858082
$arr2[17,4] = 1.067331016724282
$arr2[17,5] = 1.055407519453772
$arr2[18,0] = 1.043118311351974
$arr2[18,1] = 1.053602013496148
$arr2[18,2] = 1.052113486584878
$arr2[18,3] = 1.059189148714907
$arr2[18,4] = 1.066864909778468
$arr2[18,5] = 1.055098249205152
$arr2[19,0] = 1.042696156633778
$arr2[19,1] = 1.052595889661757
$arr2[19,2] = 1.051235314904972
$arr2[19,3] = 1.057795563628937
$arr2[19,4] = 1.065348660159681
$arr2[19,5] = 1.054090696559512
$arr2[20,0] = 1.042428943735408
$arr2[20,1] = 1.051961688546286
$arr2[20,2] = 1.050681537702125
$arr2[20,3] = 1.056918554261996
$arr2[20,4] = 1.064394623718153
$arr2[20,5] = 1.053455594805714
$arr2[21,0] = 1.04257076823532
$arr2[21,1] = 1.052298044252645
$arr2[21,2] = 1.050975261340735
$arr2[21,3] = 1.057383551683173
$arr2[21,4] = 1.064900445566142
$arr2[21,5] = 1.053792428175757
$arr2[22,0] = 1.043125358631993
$arr2[22,1] = 1.053618854202181
$arr2[22,2] = 1.052128181699739
$arr2[22,3] = 1.059212499040981
$arr2[22,4] = 1.06689031826041
$arr2[22,5] = 1.05511511382692
$arr2[23,0] = 1.043761192130012
$arr2[23,1] = 1.055144663098772
$arr2[23,2] = 1.053459039566609
$arr2[23,3] = 1.058734102895347
$arr2[23,4] = 1.069196539533424
$arr2[23,5] = 1.056643089547157

$ws.Range("B2:F25").Value = $arr1
$ws.Range("I2:N25").Value = $arr2
